# "Menage d'une fonction inutilisee" - cleanup of an unused function.
# - Adds a new documentation row (38) describing the "show-add-data" task of
#   editionServiceUrl (mirroring row 37, the "show-edit-data" task), with a
#   note that the URL is still to be filled in, highlighted in yellow.
# - Removes the now-redundant scratch rows' Service/Task content (rows 41-43)
#   which documented ajaxgetchildren / show-add-data / getfeatureinfo - the
#   show-add-data info has been folded into the new row 38 above, and the
#   remaining (ajaxgetchildren, getfeatureinfo) were unused leftovers.
# - Turns on word-wrap for D37 (the "show-edit-data" task cell).
# - Updates the view's scroll position / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Word-wrap the existing "show-edit-data" task cell.
$ws.Range("D37").WrapText = $true

# New row 38: mirrors row 37 (same file / request object / service), but for
# the "show-add-data" task, whose URL still needs to be filled in - flagged
# with a yellow highlight.
$ws.Range("A38").Value = $ws.Range("A37").Value2
$ws.Range("B38").Value = $ws.Range("B37").Value2
$ws.Range("C38").Value = $ws.Range("C37").Value2
$ws.Range("D38").Value = "show-add-data"
$ws.Range("D38").Interior.Color = $ws.Range("D36").Interior.Color
$ws.Range("E38").Value = "a compléter (l'url est renvoyée depuis le serveur)"
$ws.Range("E38").Interior.Color = 65535

# Clean out the unused leftover rows - only the style remains.
$ws.Range("C41").ClearContents()
$ws.Range("D41").ClearContents()
$ws.Range("C42").ClearContents()
$ws.Range("D42").ClearContents()
$ws.Range("C43").ClearContents()
$ws.Range("D43").ClearContents()

# Update the saved view state (scrolled/selected area).
$ws.Activate() | Out-Null
$excel.ActiveWindow.TopLeftCell = $ws.Range("B34")
$ws.Range("C41:D41").Select() | Out-Null
